# Weekly update: insert two new Pimiento price rows (Zafiro rojo / Zafiro
# verde, fecha 44488) ahead of the existing data in rows 96-97, pushing the
# rest of the table down by two rows (old row 96 -> new row 98, ...,
# old row 157 -> new row 159).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 96.
$ws.Rows.Item(96).Insert()
$ws.Rows.Item(96).Insert()

# New row 96: Zafiro rojo
$ws.Cells.Item(96,1).Value  = 7
$ws.Cells.Item(96,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(96,3).Value  = "Ñuble"
$ws.Cells.Item(96,4).Value  = 44488
$ws.Cells.Item(96,5).Value  = 16
$ws.Cells.Item(96,6).Value  = 100112002
$ws.Cells.Item(96,7).Value  = "Pimiento"
$ws.Cells.Item(96,8).Value  = "Zafiro rojo"
$ws.Cells.Item(96,9).Value  = "Primera"
$ws.Cells.Item(96,10).Value = 100
$ws.Cells.Item(96,11).Value = 43000
$ws.Cells.Item(96,12).Value = 44000
$ws.Cells.Item(96,13).Value = 43500
$ws.Cells.Item(96,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(96,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96,16).Value = 2900
$ws.Cells.Item(96,17).Value = 15
$ws.Cells.Item(96,18).Value = "Hortaliza"

# New row 97: Zafiro verde
$ws.Cells.Item(97,1).Value  = 7
$ws.Cells.Item(97,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(97,3).Value  = "Ñuble"
$ws.Cells.Item(97,4).Value  = 44488
$ws.Cells.Item(97,5).Value  = 16
$ws.Cells.Item(97,6).Value  = 100112002
$ws.Cells.Item(97,7).Value  = "Pimiento"
$ws.Cells.Item(97,8).Value  = "Zafiro verde"
$ws.Cells.Item(97,9).Value  = "Primera"
$ws.Cells.Item(97,10).Value = 100
$ws.Cells.Item(97,11).Value = 41000
$ws.Cells.Item(97,12).Value = 42000
$ws.Cells.Item(97,13).Value = 41500
$ws.Cells.Item(97,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(97,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(97,16).Value = 2767
$ws.Cells.Item(97,17).Value = 15
$ws.Cells.Item(97,18).Value = "Hortaliza"
